$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscritos")

# --- Header row ---
$ws.Range("G1").Value = "INSTITUCIÓN"
$ws.Range("G1:G10").ColumnWidth = 50

# --- Rows 2-4: existing participant (Ramirez Buendia) keeps A-F, just gains column G ---
$ws.Range("G2").Value = "ESPE"
# G3 intentionally left blank (matches source diff: no institución for that row)
$ws.Range("G4").Value = "Escuela Politécnica Nacional"

# --- Row 5-6: Chasiloa Paez Mirian Amparo ---
$ws.Range("A5").Value = "Chasiloa Páez Mirian Amparo"
# B5 identification number must stay TEXT, not be auto-converted to a number
$ws.Range("B5").NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "1709691495"
$ws.Range("C5").Value = "jennymq44@hotmail.com"
$ws.Range("D5").Value = "+593 0998179204"
$ws.Range("E5").Value = "Calderón"
$ws.Range("F5").Value = "Estudiante EPN Postgrado"
$ws.Range("G5").Value = "Escuela Politécnica Nacional"

$ws.Range("A6").Value = "Chasiloa Páez Mirian Amparo"
$ws.Range("B5").Copy($ws.Range("B6"))
$ws.Range("C6").Value = "jennymq44@hotmail.com"
$ws.Range("D6").Value = "+593 0998179204"
$ws.Range("E6").Value = "Calderón"
$ws.Range("F6").Value = "Autores - Primer Artículo"
$ws.Range("G6").Value = "EPN"

# --- Row 7: Quiguango Rivera Alisson Lorena ---
$ws.Range("A7").Value = "Quiguango Rivera Alisson Lorena"
$ws.Range("B7").NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = "1003834627"
$ws.Range("C7").Value = "jennymq44@gmail.com"
$ws.Range("D7").Value = "+593 0338138035"
$ws.Range("E7").Value = "Calderón"
$ws.Range("F7").Value = "Profesionales y Profesores Externos"
$ws.Range("G7").Value = "PUCE"

# --- Row 8: Buendía Ramírez Sebastián Daniel ---
$ws.Range("A8").Value = "Buendía Ramírez Sebastián Daniel"
$ws.Range("B2").Copy($ws.Range("B8"))
$ws.Range("C8").Value = "danielsanramirez@hotmail.com"
$ws.Range("D8").Value = "+593 987890751"
$ws.Range("E8").Value = "La Florida"
$ws.Range("F8").Value = "Autores - Primer Artículo"
$ws.Range("G8").Value = "UDLA"

# --- Rows 9-10: Ramírez Daniel ---
$ws.Range("A9").Value = "Ramírez Daniel"
$ws.Range("B2").Copy($ws.Range("B9"))
$ws.Range("C9").Value = "dsramirez1193@gmail.com"
$ws.Range("D9").Value = "+297 4846518"
$ws.Range("E9").Value = "La Vicentina"
$ws.Range("F9").Value = "Profesionales y Profesores Externos"
# G9 intentionally left blank (matches source diff)

$ws.Range("A10").Value = "Ramírez Daniel"
$ws.Range("B2").Copy($ws.Range("B10"))
$ws.Range("C10").Value = "dsramirez1193@gmail.com"
$ws.Range("D10").Value = "+297 4846518"
$ws.Range("E10").Value = "La Vicentina"
$ws.Range("F10").Value = "Autores - Primer Artículo"
$ws.Range("G10").Value = "UDLA"
